$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.09"
$ws.Range("E2").Value = "1BNBBNBBestin24h"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.521"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.382"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.472"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.075"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8020"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1417"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07420"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03258"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02992"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09257"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001667"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.250"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04716"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005738"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006267"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001047"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0004775"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.980"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.140"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1290"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04181"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007019"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003499"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1043"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008853"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005491"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6798"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03059"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Output "Applied all cell updates"